$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold, centered, bordered header style) from the
# existing header cell H1 onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2-16 for columns I (I0) and J (IF)
$data = @(
    @(6, 7),
    @(3, 5),
    @(1, 1),
    @(1, 2),
    @(4, 6),
    @(9, 9),
    @(6, 6),
    @(4, 6),
    @(8, 8),
    @(6, 7),
    @(6, 7),
    @(9, 9),
    @(5, 5),
    @(8, 8),
    @(5, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
